# Fruta / hortaliza, semanal
# Insert 3 new weekly records for "Macroferia Regional de Talca - Naranja"
# right after the existing row 255, shifting all subsequent rows down by 3
# (dimension grows from A1:T276 to A1:T279).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 256..258 (existing rows 256-276 shift to 259-279).
$ws.Rows("256:258").Insert()

# New row 256: Naranja - Lane Late - Primera
$ws.Cells.Item(256, 1).Value = 5
$ws.Cells.Item(256, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(256, 3).Value = "Maule"
$ws.Cells.Item(256, 4).Value = 44461
$ws.Cells.Item(256, 5).Value = 7
$ws.Cells.Item(256, 6).Value = "Fruta"
$ws.Cells.Item(256, 7).Value = 100102
$ws.Cells.Item(256, 8).Value = "Cítricos"
$ws.Cells.Item(256, 9).Value = 100102005
$ws.Cells.Item(256, 10).Value = "Naranja"
$ws.Cells.Item(256, 11).Value = "Lane Late"
$ws.Cells.Item(256, 12).Value = "Primera"
$ws.Cells.Item(256, 13).Value = 520
$ws.Cells.Item(256, 14).Value = 6000
$ws.Cells.Item(256, 15).Value = 6500
$ws.Cells.Item(256, 16).Value = 6308
$ws.Cells.Item(256, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(256, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(256, 19).Value = 421
$ws.Cells.Item(256, 20).Value = 15

# New row 257: Naranja - Navel Late - Primera
$ws.Cells.Item(257, 1).Value = 5
$ws.Cells.Item(257, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(257, 3).Value = "Maule"
$ws.Cells.Item(257, 4).Value = 44461
$ws.Cells.Item(257, 5).Value = 7
$ws.Cells.Item(257, 6).Value = "Fruta"
$ws.Cells.Item(257, 7).Value = 100102
$ws.Cells.Item(257, 8).Value = "Cítricos"
$ws.Cells.Item(257, 9).Value = 100102005
$ws.Cells.Item(257, 10).Value = "Naranja"
$ws.Cells.Item(257, 11).Value = "Navel Late"
$ws.Cells.Item(257, 12).Value = "Primera"
$ws.Cells.Item(257, 13).Value = 350
$ws.Cells.Item(257, 14).Value = 6000
$ws.Cells.Item(257, 15).Value = 6000
$ws.Cells.Item(257, 16).Value = 6000
$ws.Cells.Item(257, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(257, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(257, 19).Value = 400
$ws.Cells.Item(257, 20).Value = 15

# New row 258: Naranja - New Hall - Segunda
$ws.Cells.Item(258, 1).Value = 5
$ws.Cells.Item(258, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(258, 3).Value = "Maule"
$ws.Cells.Item(258, 4).Value = 44461
$ws.Cells.Item(258, 5).Value = 7
$ws.Cells.Item(258, 6).Value = "Fruta"
$ws.Cells.Item(258, 7).Value = 100102
$ws.Cells.Item(258, 8).Value = "Cítricos"
$ws.Cells.Item(258, 9).Value = 100102005
$ws.Cells.Item(258, 10).Value = "Naranja"
$ws.Cells.Item(258, 11).Value = "New Hall"
$ws.Cells.Item(258, 12).Value = "Segunda"
$ws.Cells.Item(258, 13).Value = 280
$ws.Cells.Item(258, 14).Value = 5000
$ws.Cells.Item(258, 15).Value = 5000
$ws.Cells.Item(258, 16).Value = 5000
$ws.Cells.Item(258, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(258, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(258, 19).Value = 333
$ws.Cells.Item(258, 20).Value = 15
